# Applies the changes described in the commit "Correções e Atualizações dos artefatos"
# to RT_Lista_Itens.xlsx (sheet "Lista de Itens de Trabalho").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")

# Minor tab-ratio tweak from the original edit (best effort; harmless if unsupported)
$wb.Windows.Item(1).TabRatio = 0.990

# --- Update status ("Estado Atual") values -------------------------------
# Row 11 ("Especificar os requisitos iniciais do projeto"): was "Em desenvolvimento"
$ws.Range("D11").Value = "Finalizado"

# Row 12 ("Implementar os protótipos de GUI"): was "Não Iniciado"
$ws.Range("D12").Value = "Finalizado"

# Row 12 "Horas Trabalhadas" (H12) now has a value
$ws.Range("H12").Value = 5

# Row 13 ("Preparar a apresentação" / "Não Iniciado") keeps the same
# displayed text; re-asserted here only for clarity/documentation.
$ws.Range("A13").Value = "Preparar a apresentação"
$ws.Range("D13").Value = "Não Iniciado"

# --- Add a new blank formatted row (row 16), matching the look of row 15 -
$ws.Range("B15:I15").Copy() | Out-Null
$ws.Range("B16:I16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update the selected cell in the sheet view ---------------------------
$ws.Range("D13").Select() | Out-Null
